$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B71: was stored as text "3", should be numeric 3
$ws.Range("B71").Value = 3

# Add new row 72 with the new annotation data
$ws.Range("A72").Value = "Ying Tang"
$ws.Range("B72").Value = "'3"
$ws.Range("B72").Style = "Normal"
$ws.Range("C72").Value = "无"
$ws.Range("D72").Value = "QSN"
$ws.Range("E72").Value = "MET"
$ws.Range("F72").Value = "6649e081-efd7-424b-ac96-c0694d3eea45"
$ws.Range("G72").Value = "HyydRMZC-_annotated.xlsx"
$ws.Range("H72").Value = "In particular, what is exactly Opt attack?"
